$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the missing X7/Y7 cells on row 7
$ws.Range("X7").Value = -3.5699769999999944
$ws.Range("Y7").Value = "Down"

# Add a brand new row 8 of data
$ws.Range("A8").Value = 42649.890949074077
$ws.Range("A8").NumberFormat = "m/d/yy h:mm"
$ws.Range("B8").Value = -16
$ws.Range("C8").Value = "Strong Sell"
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = "Random"
$ws.Range("Q8").Value = 49.72799223503381
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = 0.1095
$ws.Range("S8").NumberFormat = "0.00%"
$ws.Range("T8").Value = 0.0025000000000000001
$ws.Range("T8").NumberFormat = "0.00%"
$ws.Range("U8").Value = 5.95
$ws.Range("V8").Value = "N/A"
$ws.Range("W8").Value = 0
